$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ListOfValues sheet: the "Bybit_Testnet" exchange option (and its related
# strategy-settings JSON on Sheet1!K2) are being retired, so clear the cell
# that held the option while preserving its existing formatting/border style.
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("ListOfValues")
$wsList.Range("A4").ClearContents()
$wsList.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet1: main test-case table.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Policy change: the back-test window now starts 2021-01-01 instead of
# 2022-02-01 (end date 2022-12-31 is unchanged).
$ws1.Range("D2").Value = 44197

# New policy for SL/TP hitting in the same candle: instead of the old
# TP/SL percentages (1 / 0.8), both now use a fixed 0.6% with four decimal
# places of precision.
$ws1.Range("G1:H1").NumberFormat = "#,##0.0000"
$ws1.Range("G2:H2").NumberFormat = "#,##0.0000"
$ws1.Range("G2").Value = 0.006
$ws1.Range("H2").Value = 0.006

# Extend the "optional strategy settings" box down through rows 3-5 (same
# look as K2) before the worked example in K2 is cleared out.
$ws1.Rows.Item(3).RowHeight = 14.25
$ws1.Rows.Item(4).RowHeight = 14.25
$ws1.Rows.Item(5).RowHeight = 14.25

$ws1.Range("K2").Copy()
$ws1.Range("K3").PasteSpecial(-4122)
$ws1.Range("K4").PasteSpecial(-4122)
$ws1.Range("K5").PasteSpecial(-4122)

# Clear the worked UltimateScalper/Bybit_Testnet example that used to live
# in K2 - it referenced the now-removed exchange.
$ws1.Range("K2").ClearContents()

# Drop the two trailing blank template rows (15 and 16).
$ws1.Rows.Item(16).Delete()
$ws1.Rows.Item(15).Delete()

$ws1.Range("D7").Select()
